$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 'Alice Paula Di Sabatino Guimarães'
    3 = 'Ana Cláudia de Araújo Teixeira'
    4 = 'Ana Camila Oliveira Alves'
    5 = 'Adriana Costa Bacelo'
    6 = 'Anna Carolina Machado Marinho'
    7 = 'Antonio Marcos Aires Barbosa'
    8 = 'Anya Pimentel Gomes Fernandes Vieira Meyer'
    9 = 'Carla Freire Celedonio Fernandes'
    10 = 'Claudia Stutz Zubieta'
    11 = 'Clarissa Romero Teixeira'
    12 = 'Dayane Alves Costa'
    13 = 'Donat Alexander de Chapeaurouge'
    14 = 'Eduardo Ruback dos Santos'
    15 = 'Fabio Miyajima'
    16 = 'Fernando Braga Stehling Dias'
    17 = 'Fernando Ferreira Carneiro'
    18 = 'Galba Freire Moita'
    19 = 'Giovanny Augusto Camacho Antevere Mazzarotto'
    20 = 'Gilvan Pessoa Furtado'
    21 = 'Ivana Cristina de Holanda Cunha Barreto'
    22 = 'Jaime Ribeiro Filho'
    23 = 'João Hermínio Martins da Silva'
    24 = 'José Luís Passos Cordeiro'
    25 = 'Luiz Odorico Monteiro de Andrade'
    26 = 'Marcela Helena Gambim Fonseca'
    27 = 'Marcos Roberto Lourenzoni'
    28 = 'Márcio Flávio Moura de Araújo'
    29 = 'Margareth Borges Coutinho Gallo'
    30 = 'Marlos de Medeiros Chaves'
    31 = 'Maximiliano Ponte'
    32 = 'Raphael Trevizani'
    33 = 'Regis Bernardo Brandim Gomes'
    34 = 'Roberto Nicolete'
    35 = 'Roberto Wagner Júnior Freire de Freitas'
    36 = 'Sharmênia de Araújo Soares Nuto'
    37 = 'Vanira Matos Pessoa'
    38 = 'Venúcia Bruna Magalhães Pereira'
    39 = 'Fernanda Savicki de Almeida'
    40 = 'Caroline Pereira Bittencourt Passaes'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}

